# Update the Hardware-ID license list with the latest data:
#   - keep the existing "S36SNWAH859775Z" row, but move it down to row 4
#   - insert two new rows above it for "S36SNWAH859775t" with the two
#     preceding dates
#   - keep the trailing blank (but date-formatted) row after it
#   - leave the selection on C8, matching the latest saved view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Relocate the trailing empty-but-date-formatted cell from B3 to B5.
$ws.Range("B3").Copy($ws.Range("B5"))

# 2. Relocate the existing data row (row 2) down to row 4, carrying its
#    formatting (text style for A, date style for B) along with it.
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("B2").Copy($ws.Range("B4"))

# 3. Fill in the new row 2 with the latest Hardware-ID and its date.
$ws.Range("A2").Value = "S36SNWAH859775t"
$ws.Range("B2").Value = 46061

# 4. Fill in the new row 3 (B3 already carries the date style copied
#    down from the original data, so only the values need setting).
$ws.Range("A3").Value = "S36SNWAH859775t"
$ws.Range("B3").Value = 46062

# 5. Match the latest saved selection/view state.
$ws.Range("C8").Select()
